$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Shared string "2016-08-15 21:00:55" is used by both Overview!G2 and de-de!H2.
# Update both cells to the new value so the shared string stays shared.
$wsOverview.Range("G2").Value = "2016-08-15 21:01:47"
$wsDeDe.Range("H2").Value = "2016-08-15 21:01:47"

$wsZhCn.Range("H2").Value = "2016-08-15 21:01:42"
$wsZhCn.Range("K2").Value = "2016-08-15 21:01:59"

$wsDeDe.Range("K2").Value = "2016-08-15 21:02:14"
